{"js": "// Apply hybrid bold + color (#2C3E50) highlighting to quantitative impact\n// metrics (percentages, dollar amounts, large numbers) across the resume's\n// achievement / responsibility bullet points.\n//\n// Each target paragraph is located by a unique substring of its text, then\n// each metric inside that paragraph is located with a paragraph-scoped\n// search (so duplicate metrics elsewhere in the document, e.g. the same\n// \"23% to 64%\" figure quoted in the summary, are left untouched) and given\n// bold + the highlight color. Word/Office.js automatically splits the run\n// around the matched sub-range, which yields the same run structure as the\n// target OOXML (plain run, bold+colored run, plain run, ...).\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// [unique substring to find the paragraph, [ordered list of metric substrings]]\nconst targets = [\n  [\n    \"\u2022 Discovered systematic race coding errors affecting all Black and Asian-American voters\",\n    [\"23%\", \"64%\"],\n  ],\n  [\n    \"\u2022 Utilized advanced sampling methods to decrease survey margin of error\",\n    [\"\u00b14.2%\", \"\u00b12.1%\", \"71%\", \"87%\"],\n  ],\n  [\n    \"\u2022 Trigonometric algorithm for boundary estimation reduced mapping costs\",\n    [\"73.5%\", \"$4.7M\"],\n  ],\n  [\n    \"\u2022 Built real-time FEC analysis systems using Python, Pandas and PySpark\",\n    [\"$2\"],\n  ],\n  [\n    \"\u2022 Modernized legacy ETL processes by implementing dbt and PySpark workflows\",\n    [\"57%\"],\n  ],\n  [\n    \"\u2022 Platform impact: Built redistricting system serving\",\n    [\"12,847\"],\n  ],\n  [\n    \"\u2022 Revenue generation: Delivered\",\n    [\"$4.9M\"],\n  ],\n  [\n    \"\u2022 23% conversion rate improvement\",\n    [\"23%\"],\n  ],\n];\n\nfor (const [paraNeedle, metrics] of targets) {\n  const para = paragraphs.items.find((p) => p.text.indexOf(paraNeedle) !== -1);\n  if (!para) {\n    throw new Error(`Paragraph not found for: ${paraNeedle}`);\n  }\n  for (const metric of metrics) {\n    const results = para.search(metric, { matchCase: true });\n    results.load(\"items\");\n    await context.sync();\n    if (results.items.length === 0) {\n      throw new Error(`Metric \"${metric}\" not found in paragraph: ${paraNeedle}`);\n    }\n    const hit = results.items[0];\n    hit.font.bold = true;\n    hit.font.color = \"#2C3E50\";\n  }\n  await context.sync();\n}\n", "ps1": "# Apply hybrid bold + color (#2C3E50) highlighting to quantitative impact\n# metrics (percentages, dollar amounts, large numbers) across the resume's\n# achievement / responsibility bullet points.\n#\n# Each target paragraph is located by a unique substring of its text, then\n# each metric inside that paragraph is located with Find.Execute scoped to\n# a duplicate of that paragraph's Range (so duplicate metrics elsewhere in\n# the document, e.g. the same \"23% to 64%\" figure quoted in the summary,\n# are left untouched) and given bold + the highlight color. Word splits the\n# run around the matched sub-range automatically, which yields the same run\n# structure as the target OOXML (plain run, bold+colored run, plain run...).\n\n$d = $word.ActiveDocument\n\nfunction Find-ParagraphContaining($needle) {\n  foreach ($p in $d.Paragraphs) {\n    if ($p.Range.Text.Contains($needle)) {\n      return $p\n    }\n  }\n  throw \"paragraph not found: $needle\"\n}\n\nfunction Set-MetricHighlight($para, $metricText) {\n  $rng = $para.Range.Duplicate\n  $rng.Find.ClearFormatting()\n  $rng.Find.Text = $metricText\n  $rng.Find.MatchCase = $true\n  $rng.Find.MatchWildcards = $false\n  $rng.Find.Forward = $true\n  $found = $rng.Find.Execute()\n  if (-not $found) {\n    throw \"metric not found: $metricText\"\n  }\n  $rng.Font.Bold = $true\n  $rng.Font.Color = \"#2C3E50\"\n}\n\n$targets = @(\n  @{ Needle = \"\u2022 Discovered systematic race coding errors affecting all Black and Asian-American voters\"; Metrics = @(\"23%\", \"64%\") },\n  @{ Needle = \"\u2022 Utilized advanced sampling methods to decrease survey margin of error\"; Metrics = @(\"\u00b14.2%\", \"\u00b12.1%\", \"71%\", \"87%\") },\n  @{ Needle = \"\u2022 Trigonometric algorithm for boundary estimation reduced mapping costs\"; Metrics = @(\"73.5%\", \"$4.7M\") },\n  @{ Needle = \"\u2022 Built real-time FEC analysis systems using Python, Pandas and PySpark\"; Metrics = @(\"$2\") },\n  @{ Needle = \"\u2022 Modernized legacy ETL processes by implementing dbt and PySpark workflows\"; Metrics = @(\"57%\") },\n  @{ Needle = \"\u2022 Platform impact: Built redistricting system serving\"; Metrics = @(\"12,847\") },\n  @{ Needle = \"\u2022 Revenue generation: Delivered\"; Metrics = @(\"$4.9M\") },\n  @{ Needle = \"\u2022 23% conversion rate improvement\"; Metrics = @(\"23%\") }\n)\n\nforeach ($target in $targets) {\n  $para = Find-ParagraphContaining($target.Needle)\n  foreach ($metric in $target.Metrics) {\n    Set-MetricHighlight $para $metric\n  }\n}\n"}
